# Add a new "LoginLogout" worksheet between "Articles" and "RegisterUser",
# populate it with Login/Logout test data, and make it the active sheet.

$wb = $excel.ActiveWorkbook

$articles = $wb.Worksheets.Item("Articles")
$loginLogout = $wb.Worksheets.Add($null, $articles)
$loginLogout.Name = "LoginLogout"

$data = @(
    @("Key",                                      "Username",     "Password", "RememberMeCheck"),
    @("LoginWithRememberMeCheck",                  "user1@abv.bg", "user1",    $true),
    @("LogoutOfBlog",                              "user1@abv.bg", "user1",    $false),
    @("LoginWithAutocomplete",                     "user1",        "user1",    $false),
    @("LoginWithInvalidEmailAndValidPassword",     "user1",        "user1",    $false),
    @("LoginWithValidEmailAndInvalidPassword",     "user1@abv.bg", "invalid",  $false),
    @("LoginWithInvalidEmailAndInvalidPassword",   "something",    "something",$false)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    for ($c = 0; $c -lt $data[$r].Length; $c++) {
        $loginLogout.Cells.Item($r + 1, $c + 1).Value = $data[$r][$c]
    }
}

# The newly inserted sheet becomes the active / selected tab.
$loginLogout.Activate()
